$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Correct five existing adj_close text values (column G) that were
#    re-computed by the R script (tiny floating point differences).
#    Column G stores these as text, so switch it to a text number
#    format first to keep Replace() from turning the results back
#    into plain numbers.
# ------------------------------------------------------------------
$ws.Columns("G").NumberFormat = "@"

$corrections = @{
    "4.09344339370728" = "4.09344291687012"
    "4.84022045135498" = "4.84021997451782"
    "4.54387998580933" = "4.54387950897217"
    "4.7414402961731"  = "4.74143981933594"
    "4.52412414550781" = "4.52412366867065"
}

foreach ($old in $corrections.Keys) {
    $new = $corrections[$old]
    $ws.Columns("G").Replace($old, $new, 1, 1, $false, $false, $false)
}

# ------------------------------------------------------------------
# 2) Insert two new trade rows coming from the latest R script run:
#      - a new row for 2024-06-11 (50 shares @ 5.7) right before the
#        existing 2024-06-12 row, whose timestamp also gets corrected
#        from an intraday time to the plain trading date.
#      - a new row for 2024-06-13 (40 shares @ 5.95) appended at the
#        end of the sheet.
# ------------------------------------------------------------------

# Insert a new row 283, pushing the former row 283 down to 284.
$ws.Rows(283).Insert()

$ws.Cells.Item(283, 1).Value = 45454.2916666667
$ws.Cells.Item(283, 2).Value = 50
$ws.Cells.Item(283, 3).Value = 5.69999980926514
$ws.Cells.Item(283, 4).Value = 5.69999980926514
$ws.Cells.Item(283, 5).Value = 5.69999980926514
$ws.Cells.Item(283, 6).Value = 5.69999980926514
$ws.Cells.Item(283, 7).Value = "5.69999980926514"
$ws.Cells.Item(283, 8).Value = "RES.MI"

# Fix up the timestamp of the row that used to be 283 (now 284): the
# intraday time component is dropped, keeping just the trading date.
$ws.Cells.Item(284, 1).Value = 45455.2916666667

# Append the brand-new last row (285), re-using the date format from
# the row above it so the new cell keeps the same style as the rest
# of column A.
$ws.Cells.Item(285, 1).Value = 45456.6376273148
$ws.Cells.Item(284, 1).Copy()
$ws.Cells.Item(285, 1).PasteSpecial(-4122)
$ws.Cells.Item(285, 1).Value = 45456.6376273148

$ws.Cells.Item(285, 2).Value = 40
$ws.Cells.Item(285, 3).Value = 5.94999980926514
$ws.Cells.Item(285, 4).Value = 5.94999980926514
$ws.Cells.Item(285, 5).Value = 5.94999980926514
$ws.Cells.Item(285, 6).Value = 5.94999980926514
$ws.Cells.Item(285, 7).NumberFormat = "@"
$ws.Cells.Item(285, 7).Value = "5.94999980926514"
$ws.Cells.Item(285, 8).Value = "RES.MI"
